$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 20 with the new day's data
$ws.Range("A20").Value = 45964
$ws.Range("B20").Value = 674
$ws.Range("C20").Value = 41
$ws.Range("D20").Value = 633

# Update the selection to match the newly entered row
$ws.Range("A20:D20").Select()
